$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# -- Remote pairing matrix (rows 14-20) -----------------------------------
# Row 14 (Felix Schmidt)
$ws.Range("C14").Value = 240
$ws.Range("D14").Value = 60
$ws.Range("E14").Value = 860
$ws.Range("F14").Value = 480
$ws.Range("G14").Value = 120

# Row 15 (Ermal Gashi)
$ws.Range("E15").Value = 60
$ws.Range("F15").Value = 480
$ws.Range("G15").Value = 360
$ws.Range("H15").Value = 920

# Row 16 (Michael Baier)
$ws.Range("F16").Value = 600
$ws.Range("G16").Value = 450
$ws.Range("H16").Value = 360

# Row 17 (Jakob Stanta)
$ws.Range("F17").Value = 360

# Row 18 (Magdalena Hinterkörner) - overwrite the old "=E17" formula with a value
$ws.Range("E18").Value = 360

# Row 19 (Thomas Pinheiro de Souza)
$ws.Range("H19").Value = 420

# Row 20 (Florian Buchacher) - overwrite the old "=H19" formula with a value
$ws.Range("G20").Value = 600

# -- View state ------------------------------------------------------------
$ws.Range("J22").Select()
$excel.ActiveWindow.Zoom = 70
